$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix line endings in B8 and B10 (LF -> CRLF) ---
$ws.Range("B8").Value = "tst test 6`r`n"
$ws.Range("B10").Value = "new test`r`n"

# --- Remove the 4 obsolete comment rows (eferere, just test test, treestesfdsaete,
#     most recent comments) and re-insert 3 blank rows so everything below lands
#     exactly 1 row higher than before (net effect: one comment removed) ---
$ws.Rows("12:15").Delete()
$ws.Rows("12:14").Insert()

# Row 15 now needs the (re-numbered) comment that used to be two rows further down,
# with CRLF line breaks and no trailing newline.
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "jdslkfdsa er erakldf erkajfd aekrejalks dfjsalkre`r`na ereal;rkesa;skrjeas `r`nae r;alsejres"
$ws.Range("C15").Value = "'04-10-2023"

# --- Replace the former last row (was row 57: id 55 'qwewqe') with the new comment ---
$ws.Range("A57").Value = 56
$ws.Range("B57").Value = "ersdfas"
$ws.Range("C57").Value = "'04-14-2023"

# --- Append two brand-new comments, leaving row 59 empty (gap) ---
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "ewrwwerwe"
$ws.Range("C58").Value = "'2023-04-15T02:14:09.711Z"

$ws.Range("A60").Value = 58
$ws.Range("B60").Value = "test testtest`n"
$ws.Range("C60").Value = "'2023-04-15T19:51:43.709Z"
